# Apply edits to the "Analysis" worksheet of the workbook:
#  1. Wrap each HYPERLINK formula in A2:C25 with
#     IF(NOT(ISBLANK(...)), HYPERLINK(...), "") so a blank source cell in
#     Sheet1 no longer shows a dead "link" label.
#       A column -> Sheet1!E, B column -> Sheet1!F, C column -> Sheet1!G
#  2. Move the sheet's selected cell/range from F31 to E9.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws  = $wb.Worksheets.Item("Analysis")

$columns = @{
    "A" = "E"
    "B" = "F"
    "C" = "G"
}

for ($row = 2; $row -le 25; $row++) {
    foreach ($col in $columns.Keys) {
        $srcCol  = $columns[$col]
        $cellRef = "$col$row"
        $srcRef  = "Sheet1!$srcCol$row"
        $formula = "=IF(NOT(ISBLANK($srcRef)),HYPERLINK($srcRef, ""link""),"""")"
        $ws.Range($cellRef).Formula = $formula
    }
}

# Update the saved selection on the Analysis sheet to E9. Selecting a range
# on a non-active sheet activates that sheet as a side effect, so reselect
# Sheet1's original cell afterwards to leave the active tab/selection
# untouched (it was already Sheet1!A2 before this edit).
$ws.Range("E9").Select()
$ws1.Range("A2").Select()
